# Added Headless browser mode and Parallel browser testing
#
# TestOutputData sheet: remove the "SAMSUNG Galaxy S22 ultra 5G (Burgundy, 256 GB)"
# row (last row, row 9), and swap the 256GB Black/Violet product rows (3 & 4)
# into their new order. TestInputData sheet content is unchanged; only the
# active-cell selections move.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestInputData")
$ws2 = $wb.Worksheets.Item("TestOutputData")

# Remove the obsolete "Galaxy S22 ultra" product row from the output data sheet.
$ws2.Rows.Item(9).Delete()

# Re-order the two 256GB variants (Black now before... Violet now before Black).
$ws2.Range("A3").Value = "SAMSUNG Galaxy S24 Ultra 5G (Titanium Violet, 256 GB)"
$ws2.Range("A4").Value = "SAMSUNG Galaxy S24 Ultra 5G (Titanium Black, 256 GB)"

# Update the saved selections on each sheet (TestInputData stays the active tab).
$ws2.Range("C15").Select()
$ws1.Range("B2").Select()
